$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain text storage (it holds text-formatted price
# strings like "26.845.70" / "21.64" that would otherwise be auto-parsed
# as numbers when written through .Value).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.845.70"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.874.88"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "301.70"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.5353"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "0.3749"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").Value = "0.07199"
$ws.Range("E9").Value = "  -1.57%  "
$ws.Range("D10").Value = "21.64"
$ws.Range("D11").Value = "0.8911"
$ws.Range("D12").Value = "0.08189"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "1.874.76"
$ws.Range("E13").Value = "  +4.05%  "
$ws.Range("D14").Value = "93.37"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").Value = "5.313"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "14.87"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "0.000008538"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "26.875.85"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").Value = "4.995"
$ws.Range("D22").Value = "10.61"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("E23").Value = "  -1.66%  "
$ws.Range("D24").Value = "2.290"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "146.54"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "1.740"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "18.11"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("D29").Value = "4.716"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("D30").Value = "4.621"
$ws.Range("E30").Value = "  -5.34%  "
$ws.Range("D31").Value = "0.09116"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "0.8142"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").Value = "0.05018"
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").Value = "1.175"
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("D35").Value = "2.963"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("D36").Value = "0.6096"
$ws.Range("E36").Value = "  +5.81%  "
$ws.Range("D37").Value = "2.659"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("D40").Value = "1.073"
$ws.Range("E40").Value = "  -1.09%  "
$ws.Range("D41").Value = "6.617"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "8.886"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "0.5159"
$ws.Range("E43").Value = "  +4.46%  "
$ws.Range("D44").Value = "115.14"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").Value = "0.1498"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "10.01"
$ws.Range("E47").Value = "  -2.16%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "37.59"
$ws.Range("D50").Value = "0.06063"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "62.22"
$ws.Range("E51").Value = "  -3.59%  "
